# Rename Sheet2 -> assignments
$wb = $excel.ActiveWorkbook
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Name = "assignments"
